# Actualización automática de datos
# Appends the new row of data (row 14) to Sheet1, matching the source
# workbook's table of services: Fecha, Hora, Servicio, Costo, Cliente,
# EstadoPago, MetodoPago, Metodo_Pago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

$ws.Cells.Item($row, 1).Value = 45708        # A14 Fecha (date serial, keeps column's date format)
$ws.Cells.Item($row, 3).Value = "Corte Adulto"  # C14 Servicio
$ws.Cells.Item($row, 4).Value = 25              # D14 Costo
$ws.Cells.Item($row, 5).Value = "armando"       # E14 Cliente
$ws.Cells.Item($row, 7).Value = "Efectivo"      # G14 MetodoPago
